$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 63 (shifts existing rows 63.. down to 64..)
$ws.Rows("63:63").Insert()

# Populate the newly inserted row 63 with the new record's data
$ws.Cells.Item(63, 1).Value = 10
$ws.Cells.Item(63, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(63, 3).Value = "La Araucanía"
$ws.Cells.Item(63, 4).Value = 44915
$ws.Cells.Item(63, 5).Value = 9
$ws.Cells.Item(63, 6).Value = 100112031
$ws.Cells.Item(63, 7).Value = "Poroto verde"
$ws.Cells.Item(63, 8).Value = "Sin especificar"
$ws.Cells.Item(63, 9).Value = "Primera"
$ws.Cells.Item(63, 10).Value = 85
$ws.Cells.Item(63, 11).Value = 28000
$ws.Cells.Item(63, 12).Value = 28000
$ws.Cells.Item(63, 13).Value = 28000
$ws.Cells.Item(63, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(63, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(63, 16).Value = 1120
$ws.Cells.Item(63, 17).Value = 25
$ws.Cells.Item(63, 18).Value = "Hortaliza"
